# "changed researchgate to semanticscholar"
#
# Sheet1 was empty; populate it with a small table of paper references
# (title / date / authors / path_to_file) sourced from Semantic Scholar,
# with a bold, boxed, centered header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$hyphen = [char]0x2010   # U+2010 HYPHEN used inside "COVID-19"

# ---------------------------------------------------------------------
# Header row
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "title"
$ws.Range("B1").Value = "date"
$ws.Range("C1").Value = "authors"
$ws.Range("D1").Value = "path_to_file"

# Build the header format on A1 only, then stamp it onto B1:D1 via
# copy/paste-special so every header cell shares one clean cell style
# (bold font, thin box border, centered/top aligned) instead of each
# cell accumulating its own incremental style.
$a1 = $ws.Range("A1")
$a1.Font.Bold = $true
$a1.Borders.LineStyle = 1
$a1.Borders.Weight = 2
$a1.HorizontalAlignment = -4108   # xlCenter
$a1.VerticalAlignment = -4160     # xlTop
$a1.Copy()
$ws.Range("B1:D1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Data rows
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "A Study on the Impact and Response of RPA Adoption to the Customs Broker Industry in the Future"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2020"
$ws.Range("C2").Value = "Tae-In Kim; Joong-Geun Kim"
$ws.Range("D2").Value = ""

$ws.Range("A3").Value = "Impact of the COVID${hyphen}19 pandemic on surgical services: early experiences at a nominated COVID${hyphen}19 centre"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2020"
$ws.Range("C3").Value = "K. McBride; K. Brown; +3 authors C. Koh"
$ws.Range("D3").Value = "https://doi.org/10.1111/ans.15900"

$ws.Range("A4").Value = "Development of real-time reverse transcription recombinase polymerase amplification (RPA) for rapid detection of peste des petits ruminants virus in clinical samples and its comparison with real-time PCR test"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "2018"
$ws.Range("C4").Value = "Yuanli Li; Lin Li; +7 authors Zhiliang Wang"
$ws.Range("D4").Value = "https://doi.org/10.1038/s41598-018-35636-5"

$ws.Range("A5").Value = "Renoportal anastomosis in liver transplantation and its impact on patient outcomes: a systematic literature review"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "2019"
$ws.Range("C5").Value = "G. D'Amico; Ahmed Hassan; +8 authors C. Quintini"
$ws.Range("D5").Value = "https://doi.org/10.1111/tri.13368"

$ws.Range("A6").Value = "The timing of chemoradiotherapy after surgical resection and its impact on overall survival in glioblastoma."
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "2019"
$ws.Range("C6").Value = "R. Press; Sarah L. Shafer; +12 authors J. Zhong"
$ws.Range("D6").Value = "https://doi.org/10.1200/JCO.2019.37.15_SUPPL.2051"
